$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 67, shifting existing rows 67-124 down to 68-125.
$ws.Rows.Item(67).Insert(-4121)

# Populate the newly inserted row 67 with the new record.
$ws.Cells.Item(67, 1).Value = 4
$ws.Cells.Item(67, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(67, 3).Value = "Los Lagos"
$ws.Cells.Item(67, 4).Value = 44447
$ws.Cells.Item(67, 5).Value = 10
$ws.Cells.Item(67, 6).Value = 100112017
$ws.Cells.Item(67, 7).Value = "Apio"
$ws.Cells.Item(67, 8).Value = "Americana (o)"
$ws.Cells.Item(67, 9).Value = "Primera"
$ws.Cells.Item(67, 10).Value = 20
$ws.Cells.Item(67, 11).Value = 12000
$ws.Cells.Item(67, 12).Value = 12000
$ws.Cells.Item(67, 13).Value = 12000
$ws.Cells.Item(67, 14).Value = "$/docena de matas"
$ws.Cells.Item(67, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(67, 16).Value = 2000
$ws.Cells.Item(67, 17).Value = 6
$ws.Cells.Item(67, 18).Value = "Hortaliza"
